$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 178236.1217973733
$ws.Range("D2").Value = 167186.7212158808
$ws.Range("E2").Value = 160558.6133688586
$ws.Range("F2").Value = 160562.3548838807
$ws.Range("G2").Value = 155592.6031673624
$ws.Range("H2").Value = 204115.5868900973
$ws.Range("I2").Value = 204120.0642882576
$ws.Range("J2").Value = 204124.8178425392
$ws.Range("K2").Value = 204129.8645856476
$ws.Range("L2").Value = 204135.2226008278
$ws.Range("M2").Value = 204140.9110866603
$ws.Range("N2").Value = 204146.9504258519
$ws.Range("O2").Value = 204153.3622582698
$ws.Range("P2").Value = 204160.1695584813
$ws.Range("Q2").Value = 204167.3967180746
$ws.Range("R2").Value = 204175.0696330575
$ws.Range("S2").Value = 204183.2157966473

$ws.Range("C3").Value = 68236.12179737331
$ws.Range("D3").Value = 235422.8430132541
$ws.Range("E3").Value = 395981.4563821127
$ws.Range("F3").Value = 556543.8112659934
$ws.Range("G3").Value = 712136.4144333558
$ws.Range("H3").Value = 916252.0013234532
$ws.Range("I3").Value = 1120372.065611711
$ws.Range("J3").Value = 1324496.88345425
$ws.Range("K3").Value = 1528626.748039898
$ws.Range("L3").Value = 1732761.970640725
$ws.Range("M3").Value = 1936902.881727386
$ws.Range("N3").Value = 2141049.832153237
$ws.Range("O3").Value = 2345203.194411507
$ws.Range("P3").Value = 2549363.363969989
$ws.Range("Q3").Value = 2753530.760688063
$ws.Range("R3").Value = 2957705.830321121
$ws.Range("S3").Value = 3161889.046117768

$ws.Range("C5").Value = 213874.3750000001
$ws.Range("D5").Value = 213874.3750000001
$ws.Range("E5").Value = 213874.3750000001
$ws.Range("F5").Value = 213874.3750000001
$ws.Range("G5").Value = 213874.3750000001
$ws.Range("H5").Value = 213874.3750000001
$ws.Range("I5").Value = 213874.3750000001
$ws.Range("J5").Value = 213874.3750000001
$ws.Range("K5").Value = 213874.3750000001
$ws.Range("L5").Value = 213874.3750000001
$ws.Range("M5").Value = 213874.3750000001
$ws.Range("N5").Value = 213874.3750000001
$ws.Range("O5").Value = 213874.3750000001
$ws.Range("P5").Value = 213874.3750000001
$ws.Range("Q5").Value = 213874.3750000001
$ws.Range("R5").Value = 213874.3750000001
$ws.Range("S5").Value = 213874.3750000001

$ws.Range("C9").Value = 200690.0000000001
$ws.Range("D9").Value = 200690.0000000001
$ws.Range("E9").Value = 200690.0000000001
$ws.Range("F9").Value = 200690.0000000001
$ws.Range("G9").Value = 200690.0000000001
$ws.Range("H9").Value = 200690.0000000001
$ws.Range("I9").Value = 200690.0000000001
$ws.Range("J9").Value = 200690.0000000001
$ws.Range("K9").Value = 200690.0000000001
$ws.Range("L9").Value = 200690.0000000001
$ws.Range("M9").Value = 200690.0000000001
$ws.Range("N9").Value = 200690.0000000001
$ws.Range("O9").Value = 200690.0000000001
$ws.Range("P9").Value = 200690.0000000001
$ws.Range("Q9").Value = 200690.0000000001
$ws.Range("R9").Value = 200690.0000000001
$ws.Range("S9").Value = 200690.0000000001
